$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.485.94"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "3.030.27"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  +0.12%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "543.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "136.90"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.026.98"
$ws.Range("E8").Value = "  +0.83%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.494"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.45"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.87%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.148"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.450"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000222"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.11"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "3.545.04"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "62.548.53"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").Value = "3.039.13"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("E18").Value = "  -0.64%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.67"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "473.58"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.59"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.64%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.662"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.83%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "79.89"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.56"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.89%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.71"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.03%  "

$ws.Range("E30").Value = "  +0.22%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "25.73"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.51%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "54.63"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.87"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "457.87"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0807"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("E39").Value = "  +2.97%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.973.71"
$ws.Range("E40").Value = "  -5.80%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.43%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.05"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.55%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "26.97"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.14%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.251"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.85%  "

$ws.Range("E46").Value = "  -0.09%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.109"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "114.60"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").Value = "0.0₃0492"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.26"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
